$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Add Carol Nichols's reply comment on the same "Raw Identifiers" on page
#    XX cross-reference that Audrey Doyle's existing comment already covers.
# ---------------------------------------------------------------------------
$word.UserName = "Carol Nichols"
$word.UserInitials = "CN"

$text = $d.Content.Text
$anchor = $text.IndexOf([char]0x201C + "Raw Identifiers")
if ($anchor -ge 0) {
    $existing = $d.Comments.Item(1)
    $rng = $existing.Scope.Duplicate
} else {
    $rng = $d.Range(0, 0)
}
$null = $d.Comments.Add($rng, "It's in this appendix, on the 5th page.")

# ---------------------------------------------------------------------------
# 2. Clean up the tracked "capitalize first letter" insertions that were
#    originally recorded as two separate runs (single capital letter, then
#    the rest of the word) so each becomes one run inside one <w:ins>.
# ---------------------------------------------------------------------------
$word.UserName = "Audrey Doyle"
$word.UserInitials = "AD"

$count = $d.Revisions.Count
for ($i = $count; $i -ge 1; $i--) {
    $r = $d.Revisions.Item($i)
    if ($r.Type -eq 1) {
        $t = $r.Range.Text
        if ($t.Length -gt 1 -and $t -cmatch '^[A-Z][a-zA-Z ]* $') {
            $pos = $r.Range.Start
            $trackState = $d.TrackRevisions
            $d.TrackRevisions = $false
            $r.Reject()
            $d.TrackRevisions = $true
            $rng2 = $d.Range($pos, $pos)
            $rng2.InsertAfter($t)
            $d.TrackRevisions = $trackState
        }
    }
}
